$d = $word.ActiveDocument

# The page footer/boilerplate block that should be removed consists of three
# consecutive paragraphs, right after "LOB1053: Física III (Requisito fraco)":
#   1) an empty paragraph
#   2) "Ver no Jupiter Salvar em pdf Salvar em docx"
#   3) the "© 2020 ... Creative Commons Attribution" copyright line
# A trailing empty paragraph (and the page-break paragraph after it) must stay.

$count = $d.Paragraphs.Count
$startIdx = -1
$endIdx = -1

for ($i = 1; $i -le $count; $i++) {
    $text = $d.Paragraphs.Item($i).Range.Text
    if ($text -match "LOB1053") {
        $startIdx = $i + 1
    }
    if ($text -match "Creative Commons Attribution") {
        $endIdx = $i
    }
}

if ($startIdx -gt 0 -and $endIdx -ge $startIdx) {
    $rangeStart = $d.Paragraphs.Item($startIdx).Range.Start
    $rangeEnd = $d.Paragraphs.Item($endIdx).Range.End
    $r = $d.Range($rangeStart, $rangeEnd)
    $r.Delete()
}
